# Append a new test-run row (row 5) to the "Heap Report from Test" sheet,
# mirroring the layout/format of the previous row (row 4), and record the
# commit note as a new Description entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heap Report from Test")

# --- Clone row 4's number formatting down into row 5 (date + thousands- ---
# --- separated heap-size columns) so the new row matches the existing  ---
# --- table styling exactly, without inventing new style records.      ---
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)

$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("H4").Copy()
$ws.Range("H5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Fill in the new run's data ---
$ws.Range("A5").Value = 43411.458333333336   # 11/7/2018 11:00 AM
$ws.Range("B5").Value = "Laptop"
$ws.Range("C5").Value = "Release"
$ws.Range("D5").Value = "Factory_Class"
$ws.Range("F5").Value = 130007
$ws.Range("G5").Value = 130007
$ws.Range("H5").Value = 130007
$ws.Range("I5").Value = "Replaced adopt() in Code.h with call to CodeList() constructor"

# --- Reset the view back to the top-left cell ---
$ws.Activate()
[void]$ws.Range("A1").Select()
